$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the value for C60, which currently is empty, to the new note text.
$ws.Range("C60").Value = "Byte Stream Class, Reading and Writing Objects, Serializable"

# Column C has bestFit/autofit sizing; Excel recalculates the column width
# to match the new (longer) content automatically when the workbook is
# opened/saved. Set the resulting best-fit width explicitly to mirror that.
$ws.Columns.Item(3).ColumnWidth = 58.71928571428572
